$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New data column J ("North Seattle, WA") for samples 1-20 + header + avg
# ---------------------------------------------------------------------------
$ws.Range("J1").Value = "North Seattle, WA"

$northSeattle = @(12.95,13.91,11.366,12.11,12.61,13.77,13.297000000000001,12.54,13.39,13.56,11.74,11.16,12.85,13.42,13.12,13.63,12.63,12.2,13.45,12.71)
for ($i = 0; $i -lt $northSeattle.Length; $i++) {
    $row = $i + 2
    $ws.Range("J$row").Value = $northSeattle[$i]
}
$ws.Range("J22").Formula = "=AVERAGE(J2:J21)"

# ---------------------------------------------------------------------------
# 2. Remove the old one-row "Labs to Attu / localhost" summary table (K26:L27)
#    and replace it with the new 6-location summary table in F25:H30
# ---------------------------------------------------------------------------
$ws.Range("K26:L27").Clear()
$ws.Rows.Item(27).AutoFit()

$ws.Range("G25").Value = "Distance(miles)"
$ws.Range("H25").Value = "Latency"

$ws.Range("F26").Value = "Ontario, Canada"
$ws.Range("G26").Value = 1995.3
$ws.Range("H26").Value = 28.420999999999999

$ws.Range("F27").Value = "Renton, Washington"
$ws.Range("G27").Value = 15
$ws.Range("H27").Value = 13.553000000000001

$ws.Range("F28").Value = "North Seattle, Washington"
$ws.Range("G28").Value = 1.5
$ws.Range("H28").Value = 12.82

$ws.Range("F29").Value = "CS Labs"
$ws.Range("G29").Value = 0.1
$ws.Range("H29").Value = 6.6094999999999997

$ws.Range("F30").Value = "Localhost"
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 1.236

# ---------------------------------------------------------------------------
# 3. Column widths for the columns that now hold the new data / labels
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 21.166666666666668    # F  (~22)
$ws.Columns.Item(7).ColumnWidth = 13.053385416666666    # G  (~13.89)
$ws.Columns.Item(9).ColumnWidth = 17.608072916666668    # I  (~18.44)
$ws.Columns.Item(10).ColumnWidth = 17.608072916666668   # J  (~18.44)
$ws.Columns.Item(11).ColumnWidth = 20.276041666666668   # K  (~21.11)

# ---------------------------------------------------------------------------
# 4. View settings: zoom + selection
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("I37").Select()

# ---------------------------------------------------------------------------
# 5. Reposition the existing scatter chart (Chart 1)
# ---------------------------------------------------------------------------
$chart1Obj = $ws.ChartObjects().Item(1)
$chart1Obj.Left = 938.5
$chart1Obj.Top = 17.3
$chart1Obj.Width = 475.3
$chart1Obj.Height = 201.4

# ---------------------------------------------------------------------------
# 6. New bar chart (Chart 2): "Latency vs Distance"
# ---------------------------------------------------------------------------
$chart2Obj = $ws.ChartObjects().Add(705, 321, 631, 277)
$chart2 = $chart2Obj.Chart
$chart2.ChartType = 51   # xlColumnClustered
$chart2.SetSourceData($ws.Range("H26:H30"))
$series2 = $chart2.SeriesCollection().Item(1)
$series2.XValues = $ws.Range("G26:G30")
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Latency vs Distance"
$chart2.HasLegend = $false

$chart2.Axes(1).HasTitle = $true
$chart2.Axes(1).AxisTitle.Text = "Distance (miles)"
$chart2.Axes(2).HasTitle = $true
$chart2.Axes(2).AxisTitle.Text = "Latenccy(ms)"

Write-Host "edit complete"
